$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (Test Case 8): replace "enter correct format" with "use -h"
$ws.Range("E9").Value = "use -h"

# Row 10 (Test Case 9): replace "program terminates" with "use -h"
$ws.Range("E10").Value = "use -h"

# Row 11 (new Test Case 10)
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "e -f sample.txt -k "
$ws.Range("C11").Value = "This is sample file"
$ws.Range("E11").Value = "use -h"
$ws.Range("F11").Value = "FAIL"

# Update the selected/active cell to reflect the new state
$ws.Range("F12").Select()
